$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.8531859517097473
$ws.Range("B1").Value = 3.024480581283569
$ws.Range("C1").Value = 3.104901552200317
$ws.Range("D1").Value = 2.638166189193726
$ws.Range("E1").Value = 2.213926076889038
